$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) Update the first three summary cells to "0M"
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# 2) Insert 10 new rows right after row 3 (i.e. before what is currently row 4),
#    each ending up holding one value from the new per-iteration breakdown.
for ($i = 0; $i -lt 10; $i++) {
    $refRow = $t.Rows.Item(4)
    [void]$t.Rows.Add($refRow)
}

$newValues = @("22", "0.00002", "0.00007", "0.00004", "0.00000", "0.00004", "0.00004", "0.00004", "0.00081", "100.0")
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $t.Cell(4 + $i, 1).Range.Text = $newValues[$i]
}

# 3) Simplify the three multi-value rows (now shifted down by the 10 inserted rows)
#    back down to their single leading value.
$t.Cell(44,1).Range.Text = "100"
$t.Cell(45,1).Range.Text = "0"
$t.Cell(46,1).Range.Text = "75"
